$d = $word.ActiveDocument
$d.Content.Find.Execute("NS12_2", $false, $false, $false, $false, $false, $true, 1, $false, "NS12_4", 2)
